$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "tabla semaforo: add bg colors to all vbles"
# This sheet holds the threshold values (umbrales) that drive the
# traffic-light colouring; the "Positividad" column (H) thresholds are
# lowered for every risk tier so the conditional colouring kicks in
# earlier, and the last tier's threshold is reset to 0.
$ws.Range("H2").Value = 80
$ws.Range("H3").Value = 65
$ws.Range("H4").Value = 50
$ws.Range("H5").Value = 30

$ws.Range("H6").Value = 0
$ws.Range("H6").Font.Bold = $false

# Move the view back to the top-left corner and leave the selection on H7
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H7").Select()
